$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ProductLoan_Input ----
$ws1 = $wb.Worksheets.Item("ProductLoan_Input")

# Product id changes from text "kar2" to the numeric product id 485
$ws1.Range("B3").Value = 485

# Nominal interest rate default changes 12 -> 1
$ws1.Range("B11").Value = 1

# Maximum allowed outstanding balance changes 5000 -> 10000
$ws1.Range("B26").Value = 10000

# New GL-mapping rows appended after the existing data (rows 29-40)
$ws1.Range("A29").Value = "fundsource"
$ws1.Range("B29").Value = "Cash"

$ws1.Range("A30").Value = "loanprotfolio"
$ws1.Range("B30").Value = "Loan portfolio "

$ws1.Range("A31").Value = "interestreceivable"
$ws1.Range("B31").Value = "Interest Receivable "

$ws1.Range("A32").Value = "penaltiesreceivable"
$ws1.Range("B32").Value = "Penalties Receivable "

$ws1.Range("A33").Value = "transferinsuspense"
$ws1.Range("B33").Value = "Transfer in Suspence "

$ws1.Range("A34").Value = "feesreceivable"
$ws1.Range("B34").Value = "Fees Receivable"

$ws1.Range("A35").Value = "incomefrominterest"
$ws1.Range("B35").Value = "Income from interest"

$ws1.Range("A36").Value = "incomefrompenalties"
$ws1.Range("B36").Value = "Income from penalties"

$ws1.Range("A37").Value = "incomefromfees"
$ws1.Range("B37").Value = "Income from fees"

$ws1.Range("A38").Value = "incomefromrecoveryrepayments"
$ws1.Range("B38").Value = "Income from recovery repayments"

$ws1.Range("A39").Value = "loseswrittenoff"
$ws1.Range("B39").Value = "Losses Writtenoff "

$ws1.Range("A40").Value = "overpaymentliability"
$ws1.Range("B40").Value = "Overpayment Liability"

# Column B is widened to fit the new, longer GL-account descriptions
$ws1.Columns.Item(2).ColumnWidth = 55.8

# Move the visible selection down to the newly-added rows
$ws1.Range("A41").Select()

# ---- Sheet 2: ProductLoan_Output ----
$ws2 = $wb.Worksheets.Item("ProductLoan_Output")
$ws2.Range("B1").WrapText = $false
$ws2.Range("B1").Select()
